$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 25.44000000000054
$ws.Range("G2").Value = 0.000667105625100195
$ws.Range("H2").Value = 0.01216389069460296
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = 4.169463496254855
$ws.Range("L2").Value = "[1.4031065310340889, 6.935820461475622]"
$ws.Range("M2").Value = 0.003213332830893201
$ws.Range("N2").Value = 0.006426665661786402
$ws.Range("O2").Value = -2.037789829355541
$ws.Range("P2").Value = "[-2.7296320553713116, -1.34594760333977]"
$ws.Range("Q2").Value = 0.00000001294510187932474
$ws.Range("R2").Value = 0.00000002589020375864948
$ws.Range("S2").Value = 13.80441229119396
$ws.Range("T2").Value = "[12.30029461887912, 15.308529963508798]"
$ws.Range("W2").Value = 8.250810810810982
$ws.Range("X2").Value = 5.44960960960972
$ws.Range("Y2").Value = 11.05201201201224

$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 23.83000000000029
$ws.Range("G3").Value = 0.1114641026399719
$ws.Range("H3").Value = 0.2822913666013354
$ws.Range("K3").Value = 3.65377476433184
$ws.Range("L3").Value = "[-0.5960557573632119, 7.903605286026892]"
$ws.Range("M3").Value = 0.09154554344216859
$ws.Range("N3").Value = 0.09154554344216859
$ws.Range("O3").Value = 0.1320789704211913
$ws.Range("P3").Value = "[-1.446579199851156, 1.7107371406935385]"
$ws.Range("Q3").Value = 0.8690905877906339
$ws.Range("R3").Value = 0.8690905877906339
$ws.Range("S3").Value = 13.17145219252652
$ws.Range("T3").Value = "[10.751561710110131, 15.591342674942915]"
$ws.Range("W3").Value = 23.32906906906936
$ws.Range("X3").Value = 17.34175175175197
$ws.Range("Y3").Value = 29.31638638638674
